# Update the review sheet for CYRS , HSI
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Introduction sheet ---

# Ref Version bumped from 1 to 1.3
$ws1.Range("D7").Value = 1.3

# Last update date filled in as an actual date (was plain text before)
$ws1.Range("D9").Value = "9/2/2020"

# New history row: a second entry documenting this status update
$ws1.Range("B14").Value = 0.2
$ws1.Range("C14").Value = "T.Sharaby"
$ws1.Range("E14").Value = "9/2/2020"
$ws1.Range("G14").Value = "Update the status of "

# --- Cross review points sheet ---

# Mark the open points as Resolved instead of Open
$ws2.Range("H2:H6").Value = "Resolved"

# --- Window / selection state ---
$ws2.Range("F9").Select()
$ws1.Activate()
$ws1.Range("I21").Select()
